$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.Value2 = "'" + $value
    $range.Style = $style
}

# Row 2
Set-TextValue $ws.Range("D2") '27.001.69'
Set-TextValue $ws.Range("E2") '  +2.15%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.845.13'
Set-TextValue $ws.Range("E3") '  +1.89%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.04%  '

# Row 5
Set-TextValue $ws.Range("D5") '310.42'
Set-TextValue $ws.Range("E5") '  +1.26%  '

# Row 6
Set-TextValue $ws.Range("E6") '  +0.02%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.4662'
Set-TextValue $ws.Range("E7") '  +3.37%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3627'
Set-TextValue $ws.Range("E8") '  +1.13%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.07141'
Set-TextValue $ws.Range("E9") '  +0.93%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.9161'
Set-TextValue $ws.Range("E10") '  +2.50%  '

# Row 11
Set-TextValue $ws.Range("D11") '19.60'
Set-TextValue $ws.Range("E11") '  +1.12%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.07699'
Set-TextValue $ws.Range("E12") '  -1.43%  '

# Row 13
Set-TextValue $ws.Range("D13") '1.890.49'
Set-TextValue $ws.Range("E13") '  +4.58%  '

# Row 14
Set-TextValue $ws.Range("D14") '5.284'
Set-TextValue $ws.Range("E14") '  -0.07%  '

# Row 15
Set-TextValue $ws.Range("D15") '6.417'
Set-TextValue $ws.Range("E15") '  +1.74%  '

# Row 16
Set-TextValue $ws.Range("D16") '88.44'
Set-TextValue $ws.Range("E16") '  +3.90%  '

# Row 17
Set-TextValue $ws.Range("D17") '1.011'
Set-TextValue $ws.Range("E17") '  +0.10%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.000008586'
Set-TextValue $ws.Range("E18") '  +0.87%  '

# Row 19
Set-TextValue $ws.Range("E19") '  +0.00%  '

# Row 20
Set-TextValue $ws.Range("D20") '27.033.20'
Set-TextValue $ws.Range("E20") '  +2.12%  '

# Row 21
Set-TextValue $ws.Range("E21") '  +1.49%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.030'
Set-TextValue $ws.Range("E22") '  +1.16%  '

# Row 23
Set-TextValue $ws.Range("D23") '10.64'
Set-TextValue $ws.Range("E23") '  +1.20%  '

# Row 24
Set-TextValue $ws.Range("D24") '1.931'
Set-TextValue $ws.Range("E24") '  -1.48%  '

# Row 25
Set-TextValue $ws.Range("D25") '152.73'
Set-TextValue $ws.Range("E25") '  +0.46%  '

# Row 26
Set-TextValue $ws.Range("D26") '18.39'
Set-TextValue $ws.Range("E26") '  +3.32%  '

# Row 27
Set-TextValue $ws.Range("D27") '2.063'
Set-TextValue $ws.Range("E27") '  +0.36%  '

# Row 28
Set-TextValue $ws.Range("D28") '114.24'
Set-TextValue $ws.Range("E28") '  +1.68%  '

# Row 29
Set-TextValue $ws.Range("D29") '4.918'
Set-TextValue $ws.Range("E29") '  +1.37%  '

# Row 30
Set-TextValue $ws.Range("D30") '0.08860'
Set-TextValue $ws.Range("E30") '  +1.88%  '

# Row 31
Set-TextValue $ws.Range("D31") '3.187'
Set-TextValue $ws.Range("E31") '  +2.41%  '

# Row 32
Set-TextValue $ws.Range("D32") '2.854'
Set-TextValue $ws.Range("E32") '  +2.20%  '

# Row 33
Set-TextValue $ws.Range("D33") '1.174'
Set-TextValue $ws.Range("E33") '  +5.72%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.7473'
Set-TextValue $ws.Range("E34") '  +2.71%  '

# Row 35
Set-TextValue $ws.Range("D35") '4.471'
Set-TextValue $ws.Range("E35") '  +0.39%  '

# Row 36
Set-TextValue $ws.Range("D36") '1.084'
Set-TextValue $ws.Range("E36") '  +0.66%  '

# Row 37
Set-TextValue $ws.Range("B37") 'MXToken'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D37") '2.979'
Set-TextValue $ws.Range("E37") '  +2.83%  '

# Row 38
Set-TextValue $ws.Range("B38") 'VeChain'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D38") '0.01941'
Set-TextValue $ws.Range("E38") '  +0.70%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.05167'
Set-TextValue $ws.Range("E39") '  +0.92%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.5181'
Set-TextValue $ws.Range("E40") '  +1.84%  '

# Row 41
Set-TextValue $ws.Range("D41") '6.905'
Set-TextValue $ws.Range("E41") '  +2.02%  '

# Row 42
Set-TextValue $ws.Range("E42") '  -0.34%  '

# Row 43
Set-TextValue $ws.Range("D43") '8.158'
Set-TextValue $ws.Range("E43") '  +1.44%  '

# Row 44
Set-TextValue $ws.Range("D44") '10.48'
Set-TextValue $ws.Range("E44") '  +4.60%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.4702'
Set-TextValue $ws.Range("E45") '  +0.86%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +0.09%  '

# Row 47
Set-TextValue $ws.Range("D47") '100.79'
Set-TextValue $ws.Range("E47") '  +0.96%  '

# Row 48
Set-TextValue $ws.Range("D48") '1.607'
Set-TextValue $ws.Range("E48") '  +2.19%  '

# Row 49
Set-TextValue $ws.Range("D49") '64.91'
Set-TextValue $ws.Range("E49") '  +2.17%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.06047'
Set-TextValue $ws.Range("E50") '  +0.97%  '

# Row 51
Set-TextValue $ws.Range("E51") '  +0.75%  '
